# Updated symbol list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '332.41'
    'E2' = '0.81%'
    'D3' = '41.52'
    'E3' = '2.32%'
    'D4' = '5.705'
    'E4' = '-4.10%'
    'D5' = '0.08129'
    'E5' = '-0.14%'
    'D6' = '2.070'
    'E6' = '5.91%'
    'D7' = '8.749'
    'E7' = '-0.14%'
    'D8' = '4.528'
    'E8' = '-0.72%'
    'E9' = '0.57%'
    'D10' = '0.9278'
    'E10' = '-1.66%'
    'E11' = '-2.50%'
    'D12' = '0.1961'
    'E12' = '-2.36%'
    'D13' = '8.834'
    'E13' = '14.51%'
    'D14' = '0.09292'
    'E14' = '0.44%'
    'D15' = '0.03722'
    'E15' = '8.44%'
    'D16' = '0.1053'
    'E16' = '9.38%'
    'D17' = '0.001299'
    'E17' = '-1.64%'
    'D18' = '0.006218'
    'E18' = '-0.16%'
    'D19' = '3.382'
    'E19' = '0.25%'
    'E20' = '0.85%'
    'D21' = '0.1414'
    'E21' = '-1.73%'
    'D22' = '0.2605'
    'E22' = '6.52%'
    'D23' = '0.04421'
    'E23' = '-0.19%'
    'D24' = '0.001257'
    'E24' = '0.28%'
    'D25' = '0.004473'
    'E25' = '2.50%'
    'D26' = '0.0001240'
    'E26' = '4.30%'
    'D39' = '0.02922'
    'E39' = '17.15%'
    'D40' = '0.05524'
    'E40' = '4.19%'
    'D41' = '0.007734'
    'E41' = '1.68%'
    'D42' = '0.009887'
    'E42' = '10.93%'
    'E43' = '-0.63%'
    'D44' = '0.002210'
    'E44' = '7.68%'
    'D45' = '0.01100'
    'E45' = '4.02%'
    'D46' = '0.00006797'
    'E46' = '-0.32%'
    'D47' = '0.00000000750'
    'E47' = '0.08%'
    'D48' = '0.002994'
    'E48' = '3.42%'
    'E49' = '26.72%'
    'D50' = '0.00002100'
    'E50' = '0.08%'
    'D51' = '0.0002000'
    'E51' = '0.08%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

